$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.201.15"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.83%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.521.98"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -5.28%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.94"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -3.81%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.92"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.02%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.79%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.521.60"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -5.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.160"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -4.87%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.61%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -4.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.990.16"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.69%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "70.079.46"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.65%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "24.91"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.531.67"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -4.69%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.41"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -6.89%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.55"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -8.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "353.75"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -4.49%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -5.78%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.97"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -3.39%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "68.91"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -4.29%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -5.89%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -5.45%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.653.87"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -5.14%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.29%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0911"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -5.90%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "481.19"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -4.04%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.30%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.75"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -4.03%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "157.20"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.55%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.115"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.75%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.85"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.66%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -4.95%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -5.70%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -6.96%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -4.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.71"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -5.43%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -6.70%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "38.30"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "142.15"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -9.36%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.52"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -5.79%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.523"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -6.21%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -7.00%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.24%  "
